# Append a new effort-log entry (row 12) to the worksheet, mirroring the
# existing rows: a date in column A (with the same date number format as
# the rows above it), effort hours in B/C, and a task description in D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row just below the current last row (row 11).
$newRow = 12

# Column A: date value. Column A already carries the "ddd dd/mm/yyyy" date
# style at the column level, so the new cell picks it up automatically.
$ws.Cells.Item($newRow, 1).Value = "6/14/2013"

# Columns B and C: the day's effort and additional effort, in hours.
$ws.Cells.Item($newRow, 2).Value = 1.5
$ws.Cells.Item($newRow, 3).Value = 2.5

# Column D: task description (new shared string).
$ws.Cells.Item($newRow, 4).Value = "Implementation of semaphores and first, very preliminary but successfuls tests"

# Move the active selection to the new row's first cell, as in the source
# workbook.
$ws.Range("A12").Select()
